$d = $word.ActiveDocument

# 1. Update both date occurrences 06/07/2018 -> 11/07/2018
$d.Content.Find.Execute("06/07/2018", $true, $false, $false, $false, $false,
                         $true, 1, $false, "11/07/2018", 2)

# 2. Swap the two adjacent list item paragraphs:
#    "Mise en place de règles pour les fichiers d'import"            (first,  A)
#    "Vérification de la robustesse du système d'import selon la    (second, B)
#     configuration, les types d'objets et les types d'import"
# Find.Execute's Replace parameter triggers smart-quote autocorrection on the
# straight apostrophes in this text, so instead locate each run's range via
# Find and set its .Text property directly (no autocorrect there) to perform
# the swap safely, without ever creating a duplicate/ambiguous string.

$textA = "Mise en place de règles pour les fichiers d'import"
$textB = "Vérification de la robustesse du système d'import selon la configuration, les types d'objets et les types d'import"

# Locate A first (it precedes B in the document) and overwrite it with a
# unique placeholder so the two strings never collide while we work.
$placeholder = "###SWAP_PLACEHOLDER_9f3c1b###"

$rngA = $d.Content
$rngA.Find.Execute($textA) | Out-Null
$rngA.Text = $placeholder

# Now find B (still unique) and overwrite it with A's original text.
$rngB = $d.Content
$rngB.Find.Execute($textB) | Out-Null
$rngB.Text = $textA

# Finally find the placeholder and overwrite it with B's original text.
$rngC = $d.Content
$rngC.Find.Execute($placeholder) | Out-Null
$rngC.Text = $textB
